$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Ndimas: fix "max ongkir" (K column) value for the GOPAYDAY row (row 4) —
# was a placeholder 9999999, correct cap is 100000.
$ws.Range("K4").Value = 100000

# Reflect the author's final selection/active cell on that row as seen in
# the saved workbook.
$ws.Activate()
$ws.Range("K4").Select()
